$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C values (rows 1-36) from the old date-serial number to 100
$ws.Range("C1:C36").Value = 100

# Update the active selection on the sheet
$ws.Range("E7").Select()
